$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 928.5333000000001
$ws.Range("I11").Value = 928.5333000000001
$ws.Range("K11").Value = 928.5333000000001
$ws.Range("M11").Value = -788.5333000000001

# Sheet ALC, Row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 703
$ws.Range("I18").Value = 567
$ws.Range("J18").Value = 975
$ws.Range("K18").Value = 567
$ws.Range("L18").Value = 975
$ws.Range("M18").Value = -283
$ws.Range("N18").Value = -1543

# Sheet ALC, Row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I103").Value = 1502.5
$ws.Range("J103").Value = 1101.7
$ws.Range("K103").Value = 4507.5
$ws.Range("L103").Value = 3305.1
$ws.Range("M103").Value = -3921.5
$ws.Range("N103").Value = -4477.1

# Sheet ALC, Row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 13339882
$ws.Range("I106").Value = 16672352
$ws.Range("J106").Value = 10000
$ws.Range("K106").Value = 16672352
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = -16671721
$ws.Range("N106").Value = -11262

# Sheet ALC, Row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 76500
$ws.Range("I116").Value = 76500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 76500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -73058
$ws.Range("N116").ClearContents()

# Sheet ALC, Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8133.1177
$ws.Range("I132").Value = 8810.866
$ws.Range("K132").Value = 26432.598
$ws.Range("M132").Value = -23902.598

# Sheet ALC, Row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 307872.25
$ws.Range("I140").Value = 24995
$ws.Range("J140").Value = 402164.66
$ws.Range("K140").Value = 24995
$ws.Range("L140").Value = 402164.66
$ws.Range("M140").Value = -19815
$ws.Range("N140").Value = -412524.66

# Sheet ARM, Row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 12370.272
$ws.Range("I21").Value = 5178.8335
$ws.Range("J21").Value = 21000
$ws.Range("K21").Value = 5178.8335
$ws.Range("L21").Value = 21000
$ws.Range("M21").Value = -4804.8335
$ws.Range("N21").Value = -21748

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1769.7291
$ws.Range("I32").Value = 1769.7291
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1769.7291
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1482.7291
$ws.Range("N32").ClearContents()

# Sheet ARM, Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2812.8298
$ws.Range("I45").Value = 2538.4827
$ws.Range("K45").Value = 2538.4827
$ws.Range("M45").Value = -2161.4827

# Sheet ARM, Row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 56246.25
$ws.Range("I55").Value = 15000
$ws.Range("J55").Value = 69995
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 69995
$ws.Range("M55").Value = -14685
$ws.Range("N55").Value = -70625

# Sheet ARM, Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1847.2963
$ws.Range("I74").Value = 1715.5714
$ws.Range("J74").Value = 2308.3333
$ws.Range("K74").Value = 1715.5714
$ws.Range("L74").Value = 2308.3333
$ws.Range("M74").Value = -841.5714
$ws.Range("N74").Value = -4056.3333

# Sheet ARM, Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1847.2963
$ws.Range("I77").Value = 1715.5714
$ws.Range("J77").Value = 2308.3333
$ws.Range("K77").Value = 8577.857
$ws.Range("L77").Value = 11541.6665
$ws.Range("M77").Value = -4209.857
$ws.Range("N77").Value = -20277.6665

# Sheet ARM, Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5632.75
$ws.Range("J122").Value = 7580.5454
$ws.Range("L122").Value = 22741.6362
$ws.Range("N122").Value = -27641.6362

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3031.0322
$ws.Range("I132").Value = 2725.3333
$ws.Range("K132").Value = 8175.999899999999
$ws.Range("M132").Value = -5645.999899999999

# Sheet BSM, Row 109
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Sheet CRP, Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3469.4285
$ws.Range("I31").Value = 1860.88
$ws.Range("K31").Value = 1860.88
$ws.Range("M31").Value = -1565.88

# Sheet CRP, Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3469.4285
$ws.Range("I34").Value = 1860.88
$ws.Range("K34").Value = 1860.88
$ws.Range("M34").Value = -1658.88

# Sheet CRP, Row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8002.3335
$ws.Range("I86").Value = 8003.5
$ws.Range("K86").Value = 8003.5
$ws.Range("M86").Value = -6880.5

# Sheet CRP, Row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 8002.3335
$ws.Range("I89").Value = 8003.5
$ws.Range("K89").Value = 40017.5
$ws.Range("M89").Value = -34401.5

# Sheet CRP, Row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 84996
$ws.Range("J116").Value = 84996
$ws.Range("L116").Value = 84996
$ws.Range("N116").Value = -94174

# Sheet CRP, Row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1681.2
$ws.Range("I122").Value = 1708.4286
$ws.Range("K122").Value = 5125.2858
$ws.Range("M122").Value = -2675.2858

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1530.6
$ws.Range("I132").Value = 1425.6428
$ws.Range("K132").Value = 4276.928400000001
$ws.Range("M132").Value = -1746.928400000001

# Sheet CUL, Row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1231.5834
$ws.Range("I121").Value = 458.66666
$ws.Range("K121").Value = 1375.99998
$ws.Range("M121").Value = -65.99998000000005

# Sheet GSM, Row 41
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4787.75
$ws.Range("I41").Value = 4787.75
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 4787.75
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4432.75
$ws.Range("N41").ClearContents()

# Sheet GSM, Row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 77991.8
$ws.Range("J114").Value = 59239.75
$ws.Range("L114").Value = 59239.75
$ws.Range("N114").Value = -67917.75

# Sheet LTW, Row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 312.5862
$ws.Range("I55").Value = 215.93333
$ws.Range("J55").Value = 416.14285
$ws.Range("K55").Value = 215.93333
$ws.Range("L55").Value = 416.14285
$ws.Range("M55").Value = -42.93333000000001
$ws.Range("N55").Value = -762.14285

# Sheet LTW, Row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 166668320
$ws.Range("I93").Value = 1000000000
$ws.Range("J93").Value = 1976
$ws.Range("K93").Value = 1000000000
$ws.Range("L93").Value = 1976
$ws.Range("M93").Value = -999998752
$ws.Range("N93").Value = -4472

# Sheet LTW, Row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3999
$ws.Range("I100").Value = 3698.3333
$ws.Range("K100").Value = 3698.3333
$ws.Range("M100").Value = -3157.3333

# Sheet WVR, Row 12
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 6009.3335
$ws.Range("I12").Value = 20006
$ws.Range("J12").Value = 3210
$ws.Range("K12").Value = 20006
$ws.Range("L12").Value = 3210
$ws.Range("M12").Value = -19864
$ws.Range("N12").Value = -3494

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2850.5356
$ws.Range("I136").Value = 1888
$ws.Range("K136").Value = 5664
$ws.Range("M136").Value = -3114

Write-Output "done"
